$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "26.113.65"
Set-TextValue "E2" "  -0.74%  "
Set-TextValue "D3" "1.667.09"
Set-TextValue "E3" "  -1.39%  "
Set-TextValue "E4" "  -0.45%  "
Set-TextValue "D5" "209.27"
Set-TextValue "E5" "  -3.91%  "
Set-TextValue "D6" "0.5235"
Set-TextValue "E6" "  -2.02%  "
Set-TextValue "E7" "  -0.43%  "
Set-TextValue "D8" "0.2614"
Set-TextValue "E8" "  -3.94%  "
Set-TextValue "D9" "0.06334"
Set-TextValue "E9" "  -1.40%  "
Set-TextValue "D10" "21.08"
Set-TextValue "E10" "  -3.04%  "
Set-TextValue "D11" "0.07527"
Set-TextValue "E11" "  -2.34%  "
Set-TextValue "D12" "1.676.56"
Set-TextValue "E12" "  -1.06%  "
Set-TextValue "D13" "4.425"
Set-TextValue "E13" "  -2.35%  "
Set-TextValue "D14" "0.5488"
Set-TextValue "E14" "  -5.36%  "
Set-TextValue "D15" "66.34"
Set-TextValue "E15" "  -0.84%  "
Set-TextValue "D16" "0.000007972"
Set-TextValue "E16" "  -4.74%  "
Set-TextValue "D17" "26.133.09"
Set-TextValue "E17" "  -0.70%  "
Set-TextValue "E18" "  -0.47%  "
Set-TextValue "D19" "4.706"
Set-TextValue "E19" "  -4.13%  "
Set-TextValue "D20" "186.63"
Set-TextValue "E20" "  -3.45%  "
Set-TextValue "D21" "10.25"
Set-TextValue "E21" "  -5.68%  "
Set-TextValue "D22" "6.177"
Set-TextValue "E22" "  -1.51%  "
Set-TextValue "D23" "1.003"
Set-TextValue "E23" "  -0.44%  "
Set-TextValue "D24" "149.77"
Set-TextValue "E24" "  +0.54%  "
Set-TextValue "D25" "0.1243"
Set-TextValue "E25" "  -3.40%  "
Set-TextValue "D26" "7.473"
Set-TextValue "E26" "  -5.08%  "
Set-TextValue "D27" "15.88"
Set-TextValue "E27" "  +0.05%  "
Set-TextValue "D28" "0.06368"
Set-TextValue "E28" "  +4.10%  "
Set-TextValue "D29" "1.345"
Set-TextValue "E29" "  -2.71%  "
Set-TextValue "D30" "1.272"
Set-TextValue "E30" "  -4.14%  "
Set-TextValue "D31" "3.485"
Set-TextValue "E31" "  -3.26%  "
Set-TextValue "D32" "3.408"
Set-TextValue "E32" "  -4.90%  "
Set-TextValue "D33" "1.637"
Set-TextValue "E33" "  -3.17%  "
Set-TextValue "E34" "  -3.08%  "
Set-TextValue "E35" "  -0.83%  "
Set-TextValue "D36" "0.5995"
Set-TextValue "E36" "  -3.25%  "
Set-TextValue "E37" "  -0.73%  "
Set-TextValue "B38" "Maker"
Set-TextValue "C38" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D38" "1.107.83"
Set-TextValue "E38" "  -0.26%  "
Set-TextValue "B39" "FraxShare"
Set-TextValue "C39" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D39" "6.106"
Set-TextValue "E39" "  -2.02%  "
Set-TextValue "D40" "0.01612"
Set-TextValue "E40" "  -1.68%  "
Set-TextValue "D41" "0.8673"
Set-TextValue "E41" "  -1.05%  "
Set-TextValue "E42" "  -0.90%  "
Set-TextValue "D43" "99.93"
Set-TextValue "E43" "  -0.98%  "
Set-TextValue "D44" "1.817.82"
Set-TextValue "E44" "  -1.27%  "
Set-TextValue "E45" "  +2.22%  "
Set-TextValue "D46" "55.22"
Set-TextValue "E46" "  -4.64%  "
Set-TextValue "D47" "0.9984"
Set-TextValue "E47" "  -0.99%  "
Set-TextValue "D48" "8.045"
Set-TextValue "E48" "  -0.94%  "
Set-TextValue "D49" "0.05227"
Set-TextValue "E49" "  -1.13%  "
Set-TextValue "E50" "  -1.11%  "
Set-TextValue "D51" "5.919"
Set-TextValue "E51" "  -2.31%  "
